$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D1").NumberFormat = "General"
$d1 = $ws.Range("D1")
Write-Output $d1.NumberFormat
